$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Euramet")

# Rows 7-8: update existing values (E,F,G,H)
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = -157.0411782090434
$ws.Range("G7").Value = -1.1772
$ws.Range("H7").Value = 1

$ws.Range("E8").Value = 0
$ws.Range("F8").Value = -156.5579745837849
$ws.Range("G8").Value = -1.1772
$ws.Range("H8").Value = 1

# Rows 9-19: fill previously empty cells (D,E,F,G,H)
$ws.Range("D9").Value = 837
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = -157.0411782090434
$ws.Range("G9").Value = -1.1772
$ws.Range("H9").Value = 1

$ws.Range("D10").Value = 837
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = -157.0411782090434
$ws.Range("G10").Value = -1.1772
$ws.Range("H10").Value = 1

$ws.Range("D11").Value = 837
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = -157.0411782090434
$ws.Range("G11").Value = -1.1772
$ws.Range("H11").Value = 1

$ws.Range("D12").Value = 837
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = -157.0411782090434
$ws.Range("G12").Value = -1.1772
$ws.Range("H12").Value = 1

$ws.Range("D13").Value = 837
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = -156.5579745837849
$ws.Range("G13").Value = -1.1772
$ws.Range("H13").Value = 1

$ws.Range("D14").Value = 837
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = -157.0411782090434
$ws.Range("G14").Value = -1.1772
$ws.Range("H14").Value = 1

$ws.Range("D15").Value = 837
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = -157.0411782090434
$ws.Range("G15").Value = -1.1772
$ws.Range("H15").Value = 1

$ws.Range("D16").Value = 837
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = -156.5579745837849
$ws.Range("G16").Value = -1.1772
$ws.Range("H16").Value = 1

$ws.Range("D17").Value = 837
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = -157.0411782090434
$ws.Range("G17").Value = -1.1772
$ws.Range("H17").Value = 1

$ws.Range("D18").Value = 837
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = -156.5579745837849
$ws.Range("G18").Value = -1.1772
$ws.Range("H18").Value = 1

$ws.Range("D19").Value = 837
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = -157.0411782090434
$ws.Range("G19").Value = -1.1772
$ws.Range("H19").Value = 1

# Rows 29-37: clear previously filled cells (D,E,F,G,H)
$ws.Range("D29:H29").ClearContents()
$ws.Range("D30:H30").ClearContents()
$ws.Range("D31:H31").ClearContents()
$ws.Range("D32:H32").ClearContents()
$ws.Range("D33:H33").ClearContents()
$ws.Range("D34:H34").ClearContents()
$ws.Range("D35:H35").ClearContents()
$ws.Range("D36:H36").ClearContents()
$ws.Range("D37:H37").ClearContents()
